$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at row 9 for the "Google Scholar" profile field. This
#    shifts every row from 9 downward by one (rows 9-36 -> rows 10-37).
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()

# Label cell (A9) should look like the other label cells in this block
# (A4:A10 use the same style) - copy formatting from A10 (the "Bio" label,
# which already carries that style) then set the text.
$ws.Range("A10").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Google Scholar"

# Value cell (B9) holds the Google Scholar profile URL (with the bidi
# control characters copied from the source data) and gets its own new
# left/top aligned style.
$ws.Range("B9").Value = [char]0x202A + "https://scholar.google.ca/citations?user=ghYLsSAAAAAJ&hl=en&oi=ao" + [char]0x202C
$ws.Range("B9").HorizontalAlignment = -4131
$ws.Range("B9").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 2. The row insertion above shifts the "Great Lakes Forestry Centre" bio
#    hyperlink from B23 down to B24, but this runtime does not auto-update
#    the worksheet's Hyperlinks collection when rows move. Rebuild the
#    hyperlinks so they point at the correct (shifted) cells while keeping
#    each cell's existing text/style intact.
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("Z4").PasteSpecial(-4122)

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B5"), "https://twitter.com/wet_erik?lang=en") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:erik.emilson@canada.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/ErikEmilson") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B24"), "https://www.nrcan.gc.ca/science-data/research-centres-labs/forestry-research-centres/great-lakes-forestry-centre/13459") | Out-Null

$ws.Range("Z1").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("Z4").Copy()
$ws.Range("B24").PasteSpecial(-4122)

$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------------
# 3. Touch the page setup (orientation) - matches the print-setup change
#    captured alongside the content edits.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Restore the view (selection/scroll position) to match the saved state.
# ---------------------------------------------------------------------------
$ws.Range("B9").Select()
